# Update cryptos list: price (D) and 1h-volume-change (E) columns,
# plus a few coin-identity (B/C) swaps, per the Nov 12 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.180.94'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.055.37'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.06%  '

$ws.Range('E4').Value = '  +0.23%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.99%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.666'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.43%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.17'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -6.33%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.387'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0790'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.04%  '

$ws.Range('E11').Value = '  +0.89%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.92'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.31%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.353.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.841'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.80'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.24%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.058.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.80%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +18.59%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.180.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.05%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '75.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.35%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0904'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.41%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.95%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.72%  '

$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.49'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.35%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.03%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.49%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.47%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.38%  '

$ws.Range('E29').Value = '  -0.52%  '

$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.76%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.75%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0627'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.03%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.57'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.95%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0910'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.38%  '

$ws.Range('E35').Value = '  +0.14%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.29'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.77%  '

$ws.Range('E37').Value = '  -0.61%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.107'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.24%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.35'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.82%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.08'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.55%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0223'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.60%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.54'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.56%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.15%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.99'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.79%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.58%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.288.64'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.66%  '

$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.66%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.82%  '

$ws.Range('B50').Value = 'FTXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -14.76%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.246.12'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.32%  '
